$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 4 blank rows right before the totals row (old row 15). ---
#     This pushes the old "totals" row to 19 and the footer row to 20.
for ($i = 0; $i -lt 4; $i++) {
  $ws.Rows.Item(15).Insert()
}

# --- 2. Give the 4 freshly inserted rows (15-18) the same look as the other data rows. ---
$ws.Range("A13:N13").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A17:N17").PasteSpecial(-4122)

$ws.Range("A14:N14").Copy()
$ws.Range("A16:N16").PasteSpecial(-4122)
$ws.Range("A18:N18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Recreate the B:G / H:K / L:M merges for the new rows
foreach ($r in 15..18) {
  $ws.Range("B$r`:G$r").Merge()
  $ws.Range("H$r`:K$r").Merge()
  $ws.Range("L$r`:M$r").Merge()
}

# Row heights matching the rest of the (alternating) data rows
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 25.5
$ws.Rows.Item(19).RowHeight = 25.5

# --- 3. Re-write the whole (now 15-row) product table in alphabetical order, ---
#        inserting the 4 new products at their correct sorted position.
$rows = @(
  @{ Name = "ADWIFLAM 75MG/3ML 6 AMP.";                    Qty = "1:6";    Price = 12;      Deals = "0:0" },
  @{ Name = "CETAL 250MG/5ML 60ML SUSP";                   Qty = "17:0";   Price = 31;      Deals = "1:0" },
  @{ Name = "DURICEF 500MG/5ML SUSP. 60ML";                Qty = "1:0";    Price = 78;      Deals = "1:0" },
  @{ Name = "INJECTMOL 1 GM/100ML VIAL FOR I.V. INF.";     Qty = "6:0";    Price = 67;      Deals = "1:0" },
  @{ Name = "LEVANIC 500MG 7 F.C. TAB.";                   Qty = "0:0";    Price = 92;      Deals = "1:0" },
  @{ Name = "MEGAFEN-N 100MG/5ML SUSP. 120 ML";            Qty = "2:0";    Price = 35;      Deals = "1:0" },
  @{ Name = "ORS 10 SACHET";                               Qty = "4:1";    Price = 4;       Deals = "0:0" },
  @{ Name = "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML";   Qty = "5:0";    Price = 24;      Deals = "1:0" },
  @{ Name = "PANADOL ADVANCE 500 MG 48 TABLETS";           Qty = "3:3";    Price = 23;      Deals = "0:0" },
  @{ Name = "WATER FOR INJECTION AMP. 5 ML";                Qty = "7791:0"; Price = 5;       Deals = "2:0" },
  @{ Name = "جهاز محلول ";                                  Qty = "15:0";   Price = 20;      Deals = "1:0" },
  @{ Name = "حفاضات كبار سن ماكسويل 63ق";                    Qty = "0:27";   Price = 11.67;   Deals = "0:0" },
  @{ Name = "سرنجات 3 سم";                                  Qty = "-1:0";   Price = 2;       Deals = "1:0" },
  @{ Name = "صابون ديتول العنايه بالبشره";                    Qty = "15:0";   Price = 133.65;  Deals = "3:0" },
  @{ Name = "كالونا ";                                       Qty = "-1:0";   Price = 15;      Deals = "1:0" }
)

$r = 4
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = ($r - 3)          # column A: running index
  $ws.Range("B$r").Value = $row.Name              # column B: product name
  $ws.Range("H$r").Value = $row.Qty                # column H: balance
  $ws.Range("L$r").Value = $row.Price              # column L: price
  $ws.Range("N$r").Value = $row.Deals              # column N: deals count
  $r += 1
}

# --- 4. Refresh the grand-total cell (sum of the price column). ---
$ws.Range("K19").Value = 553.32

Write-Host "edit applied"
